$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value of 45190 for every
# data row (rows 2-407). Update it to 45192 (two days later) everywhere.
$range = $ws.Range("C2:C407")
$range.Value = 45192
